$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60.83499908447266
$ws.Range("C2").Value = 45.14500045776367
$ws.Range("D2").Value = 52.18639355018491
